$wb = $excel.ActiveWorkbook

# The "FindCarTest" worksheet holds the car-brand / browser / run-mode table
$ws = $wb.Worksheets.Item("FindCarTest")

# Update the data rows: normalize browserType to "chrome" and runmode to "Y"
# for every row (removing the now-unused "firefox" / "N" shared strings).
$ws.Range("A2").Value = "bmw"
$ws.Range("B2").Value = "chrome"
$ws.Range("C2").Value = "Y"

$ws.Range("A3").Value = "mg"
$ws.Range("B3").Value = "chrome"
$ws.Range("C3").Value = "Y"

$ws.Range("A4").Value = "toyota"
$ws.Range("B4").Value = "chrome"
$ws.Range("C4").Value = "Y"

# Update the saved cursor/selection position on this sheet
$ws.Activate()
$ws.Range("E20").Select()
